$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R: "reference_no_query" header + its lookup SQL value,
# mirroring the existing *_query columns (K..Q) pattern.
$ws.Range("R1").Value = "reference_no_query"
$ws.Range("R2").Value = "SELECT DT.LEAD_FIELD5 FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '"

# Give the new column a bestFit-style width similar to its neighbours.
$ws.Columns("R").ColumnWidth = 76.3

# Move the selection to the newly added cell, scrolled into view.
$ws.Range("O2").Select() | Out-Null
